# HACKATHON.pptx - slide 1 content placeholder: fix "MUTLI-DB ADAPTER" quote.
#
# The quoted line currently reads:  " MUTLI-DB ADAPTER "
# and is built out of 4 runs:
#   1) "             " (leading spaces, red Arial)
#   2) "\u201C "        (opening curly quote, grey Arial)
#   3) "MUTLI-DB ADAPTER " (grey Algerian - contains the typo)
#   4) "\u201D"          (closing curly quote, grey Arial)
#
# The fix splits run 3 into its own "MULTI-DB " run (typo corrected) and
# an "ADAPTER " run, keeping all existing character formatting (font,
# size, bold, color) intact.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(6)             # "Content Placeholder 8"
$tf = $sh.TextFrame
$tr = $tf.TextRange

$quoteLine = $tr.Paragraphs(2)
$adapterRun = $quoteLine.Runs(3)    # "MUTLI-DB ADAPTER "

# Re-point at the first 9 characters ("MUTLI-DB ") of that run and retype
# them with the corrected spelling; the engine splits the run in two,
# leaving "ADAPTER " behind in the original run while the new leading
# chunk becomes its own run carrying the same rPr (Algerian typeface etc).
$typo = $tr.Characters($adapterRun.Start, 9)
$typo.Text = "MULTI-DB "
